$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: write a "text that looks like a date" cell without Excel
# auto-converting it to a date serial number. Force the cell to Text
# format, assign the literal string, then restore the original
# (YYYY.MM.DD) number format used by the date column.
function Set-TextDate($addr, $text) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.NumberFormat = "YYYY.MM.DD"
}

# --- Column B (Ф.И.О.) is unchanged for every existing row; only the
# Должность (C) and Дата поступления (D) columns are updated.

# Row 5 - Иванов И.М. stays Зам. директора, date changes
Set-TextDate "D5" "02.08.2010"

# Row 6 - Коробова П.Н becomes Менеджер
$ws.Range("C6").Value = "Менеджер"
Set-TextDate "D6" "10.19.2005"

# Row 7 - Морозов И.Р. (was listed lower) becomes Водитель
$ws.Range("C7").Value = "Водитель"
Set-TextDate "D7" "04.04.2008"

# Row 8 - Ромашова П.Т. becomes Секетарь
$ws.Range("C8").Value = "Секетарь"
Set-TextDate "D8" "07.14.2008"

# Row 9 - Петров Г.Т. becomes Бухгалтер
$ws.Range("C9").Value = "Бухгалтер"
Set-TextDate "D9" "08.29.2011"

# Row 10 - Смирнов С.И. becomes Директор
$ws.Range("C10").Value = "Директор"
Set-TextDate "D10" "06.05.2004"

# Row 11 - Соколова О.С. becomes Зам. директора
$ws.Range("C11").Value = "Зам. директора"
Set-TextDate "D11" "09.19.2003"

# --- Widen column B to fit the new long label text below the table
$ws.Range("B1").ColumnWidth = 28 - 0.8333333333333334

# --- Summary block (rows 14-17). Give the new rows the same
# border+fill style already used by the table header row (row 4),
# then fill in labels/values.
$ws.Range("A4").Copy()
$ws.Range("B14:C17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B14").Value = "Курс доллара: "
$ws.Range("C14").Value = 41.3

$ws.Range("B15").Value = "Средняя зарплата, руб:"

$ws.Range("B16").Value = "Максимальная зарплата, руб:"

$ws.Range("B17").Value = "Минимальная зарплата, руб:"
